$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the four values in row 24 (C24:F24) from 2 to 5
$rng = $ws.Range("C24:F24")
$rng.Value = 5

# Remove the green highlight fill on these cells (they keep their existing
# border/font/alignment, matching the "unhighlighted" style used elsewhere
# on the sheet, e.g. C25/D25/G24 which have no interior fill)
$rng.Interior.Pattern = -4142
$rng.Interior.PatternColorIndex = -4105

# Move/record the active selection to C24, as in the authored change
$ws.Range("C24").Select()
